# Update "想去人数" (want-to-go count) figures in both the "展览" sheet
# and the "全部类型" sheet to reflect the latest generated numbers.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - rows 2-4 in column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 2198
$wsExhibit.Range("F3").Value = 918
$wsExhibit.Range("F4").Value = 1698

# Sheet "全部类型" (All types) - matching rows in column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 2198
$wsAll.Range("F5").Value = 918
$wsAll.Range("F6").Value = 1698
